$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.362.68'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.842.39'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9987'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6307'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07446'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '25.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2893'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.35%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.850.63'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.969'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6761'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001024'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '81.55'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.258'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.371.72'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9999'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.367'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '157.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.497'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1350'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06931'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.463'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.481'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.051'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.048'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.826'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.140'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6981'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.584'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01844'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.816'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.821'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.50%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.233.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9347'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9998'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.994.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000119'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.029'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.708'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.931'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1140'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3916'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.49%  '
